$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Fix formatting on existing rows 7 and 9 so they match the "data row" style
# (border + vertical-center + wrap, no pattern fill) used by rows 6/8 ---
[void]$ws.Range("B6").Copy()
[void]$ws.Range("B7").PasteSpecial(-4122)   # xlPasteFormats
[void]$ws.Range("A9").PasteSpecial(-4122)   # xlPasteFormats
$excel.CutCopyMode = $false

# --- Update participant role text for row 9 (was "Supervisor") ---
$ws.Range("B9").Value = "Asesor Comercial"

# --- Append new participants (rows 10-13) ---
$ws.Range("A10").Value = 5
$ws.Range("B10").Value = "Empleado de recursos humanos"

$ws.Range("A11").Value = 6
$ws.Range("B11").Value = "Director logístico"

$ws.Range("A12").Value = 7
$ws.Range("B12").Value = "Asesor de imagen de marca"

$ws.Range("A13").Value = 8
$ws.Range("B13").Value = "Gerente general"

# Give the new rows the same formatting as the rows above them
[void]$ws.Range("A9:B9").Copy()
[void]$ws.Range("A10:B13").PasteSpecial(-4122)   # xlPasteFormats
$excel.CutCopyMode = $false

# --- Match the saved selection state ---
[void]$ws.Range("B9").Select()
